$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows.Item(3).EntireRow.Delete() | Out-Null
$ws1.Rows.Item(3).EntireRow.Delete() | Out-Null

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Rows.Item(3).EntireRow.Delete() | Out-Null
